$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (F column: 想去人数 / interest count) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 124
$ws1.Range("F4").Value = 709
$ws1.Range("F5").Value = 130
$ws1.Range("F7").Value = 43
$ws1.Range("F8").Value = 2674
$ws1.Range("F9").Value = 1643
$ws1.Range("F10").Value = 1691
$ws1.Range("F11").Value = 316
$ws1.Range("F12").Value = 276
$ws1.Range("F13").Value = 697
$ws1.Range("F14").Value = 853
$ws1.Range("F15").Value = 137
$ws1.Range("F16").Value = 352
$ws1.Range("F17").Value = 1105
$ws1.Range("F20").Value = 533
$ws1.Range("F21").Value = 5988
$ws1.Range("F22").Value = 240
$ws1.Range("F23").Value = 1190
$ws1.Range("F24").Value = 126
$ws1.Range("F26").Value = 151
$ws1.Range("F27").Value = 279
$ws1.Range("F28").Value = 240
$ws1.Range("F30").Value = 1072
$ws1.Range("F31").Value = 861
$ws1.Range("F33").Value = 76
$ws1.Range("F35").Value = 438
$ws1.Range("F36").Value = 1252
$ws1.Range("F37").Value = 150
$ws1.Range("F38").Value = 130
$ws1.Range("F41").Value = 149

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "已停售"
$ws2.Range("F3").Value = 436

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "已停售"
$ws4.Range("F3").Value = 436
$ws4.Range("F4").Value = 124
$ws4.Range("F6").Value = 709
$ws4.Range("F7").Value = 130
$ws4.Range("F12").Value = 43
$ws4.Range("F13").Value = 2674
$ws4.Range("F14").Value = 1643
$ws4.Range("F15").Value = 1691
$ws4.Range("F16").Value = 316
$ws4.Range("F17").Value = 276
$ws4.Range("F18").Value = 697
$ws4.Range("F20").Value = 853
$ws4.Range("F21").Value = 137
$ws4.Range("F22").Value = 352
$ws4.Range("F23").Value = 1105
$ws4.Range("F25").Value = 533
$ws4.Range("F26").Value = 5988
$ws4.Range("F27").Value = 240
$ws4.Range("F28").Value = 1190
$ws4.Range("F29").Value = 126
$ws4.Range("F31").Value = 151
$ws4.Range("F32").Value = 279
$ws4.Range("F33").Value = 240
$ws4.Range("F35").Value = 1072
$ws4.Range("F36").Value = 861
$ws4.Range("F38").Value = 76
$ws4.Range("F40").Value = 438
$ws4.Range("F41").Value = 1252
$ws4.Range("F42").Value = 151
$ws4.Range("F43").Value = 130
$ws4.Range("F46").Value = 149
